# Auto-generated Excel COM-interop script
# Applies the scheduled-runner market-data refresh to Raiden_Profits.xlsx
# (updates currentAveragePrice*/LevePrice*/LeveProfit* columns H:N per leve row)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 3431.0476
$ws.Cells.Item(28, 9).Value = 3244.3076
$ws.Cells.Item(28, 10).Value = 3734.5
$ws.Cells.Item(28, 11).Value = 3244.3076
$ws.Cells.Item(28, 12).Value = 3734.5
$ws.Cells.Item(28, 13).Value = -2759.3076
$ws.Cells.Item(28, 14).Value = -4704.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 137.28572
$ws.Cells.Item(80, 9).Value = 182.7
$ws.Cells.Item(80, 10).Value = 23.75
$ws.Cells.Item(80, 11).Value = 548.0999999999999
$ws.Cells.Item(80, 12).Value = 71.25
$ws.Cells.Item(80, 13).Value = 449.9000000000001
$ws.Cells.Item(80, 14).Value = -2067.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(83, 8).Value = 137.28572
$ws.Cells.Item(83, 9).Value = 182.7
$ws.Cells.Item(83, 10).Value = 23.75
$ws.Cells.Item(83, 11).Value = 1644.3
$ws.Cells.Item(83, 12).Value = 213.75
$ws.Cells.Item(83, 13).Value = 3347.7
$ws.Cells.Item(83, 14).Value = -10197.75

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(106, 8).Value = 2344.65
$ws.Cells.Item(106, 9).Value = 2283.842
$ws.Cells.Item(106, 11).Value = 2283.842
$ws.Cells.Item(106, 13).Value = -1652.842

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(125, 8).Value = 818.85
$ws.Cells.Item(125, 9).Value = 376.85715
$ws.Cells.Item(125, 10).Value = 1056.8462
$ws.Cells.Item(125, 11).Value = 3391.71435
$ws.Cells.Item(125, 12).Value = 9511.6158
$ws.Cells.Item(125, 13).Value = -931.7143499999997
$ws.Cells.Item(125, 14).Value = -14431.6158

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(132, 8).Value = 1456.7368
$ws.Cells.Item(132, 9).Value = 1324.4857
$ws.Cells.Item(132, 11).Value = 3973.4571
$ws.Cells.Item(132, 13).Value = -1443.4571

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3238.3635
$ws.Cells.Item(32, 9).Value = 2094.1904
$ws.Cells.Item(32, 11).Value = 2094.1904
$ws.Cells.Item(32, 13).Value = -1807.1904

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(45, 8).Value = 3922
$ws.Cells.Item(45, 10).Value = 4007
$ws.Cells.Item(45, 12).Value = 4007
$ws.Cells.Item(45, 14).Value = -4761

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1217.3334
$ws.Cells.Item(20, 10).Value = 1533.6666
$ws.Cells.Item(20, 12).Value = 1533.6666
$ws.Cells.Item(20, 14).Value = -2027.6666

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(94, 8).Value = 4176.8667
$ws.Cells.Item(94, 10).Value = 16499.666
$ws.Cells.Item(94, 12).Value = 16499.666
$ws.Cells.Item(94, 14).Value = -17401.666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 34470.082
$ws.Cells.Item(28, 10).Value = 34470.082
$ws.Cells.Item(28, 12).Value = 34470.082
$ws.Cells.Item(28, 14).Value = -34960.082

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(62, 8).Value = 11720.728
$ws.Cells.Item(62, 9).Value = 12350.944
$ws.Cells.Item(62, 10).Value = 8884.75
$ws.Cells.Item(62, 11).Value = 12350.944
$ws.Cells.Item(62, 12).Value = 8884.75
$ws.Cells.Item(62, 13).Value = -11726.944
$ws.Cells.Item(62, 14).Value = -10132.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(65, 8).Value = 11720.728
$ws.Cells.Item(65, 9).Value = 12350.944
$ws.Cells.Item(65, 10).Value = 8884.75
$ws.Cells.Item(65, 11).Value = 61754.72
$ws.Cells.Item(65, 12).Value = 44423.75
$ws.Cells.Item(65, 13).Value = -58634.72
$ws.Cells.Item(65, 14).Value = -50663.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 11872806
$ws.Cells.Item(4, 9).Value = 1703735.8
$ws.Cells.Item(4, 11).Value = 5111207.4
$ws.Cells.Item(4, 13).Value = -5111095.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 2198
$ws.Cells.Item(50, 9).Value = 3084.25
$ws.Cells.Item(50, 10).Value = 1016.3333
$ws.Cells.Item(50, 11).Value = 9252.75
$ws.Cells.Item(50, 12).Value = 3048.9999
$ws.Cells.Item(50, 13).Value = -8771.75
$ws.Cells.Item(50, 14).Value = -4010.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(53, 8).Value = 2198
$ws.Cells.Item(53, 9).Value = 3084.25
$ws.Cells.Item(53, 10).Value = 1016.3333
$ws.Cells.Item(53, 11).Value = 9252.75
$ws.Cells.Item(53, 12).Value = 3048.9999
$ws.Cells.Item(53, 13).Value = -8771.75
$ws.Cells.Item(53, 14).Value = -4010.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(69, 8).Value = 9998
$ws.Cells.Item(69, 9).Value = 9998
$ws.Cells.Item(69, 11).Value = 29994
$ws.Cells.Item(69, 13).Value = -29183

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(72, 8).Value = 9998
$ws.Cells.Item(72, 9).Value = 9998
$ws.Cells.Item(72, 11).Value = 89982
$ws.Cells.Item(72, 13).Value = -85926

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 10706.75
$ws.Cells.Item(80, 10).Value = 16224.75
$ws.Cells.Item(80, 12).Value = 16224.75
$ws.Cells.Item(80, 14).Value = -18220.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(83, 8).Value = 10706.75
$ws.Cells.Item(83, 10).Value = 16224.75
$ws.Cells.Item(83, 12).Value = 81123.75
$ws.Cells.Item(83, 14).Value = -91107.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(107, 8).Value = 514.8095
$ws.Cells.Item(107, 9).Value = 321
$ws.Cells.Item(107, 10).Value = 560.41174
$ws.Cells.Item(107, 11).Value = 321
$ws.Cells.Item(107, 12).Value = 560.41174
$ws.Cells.Item(107, 13).Value = 1599
$ws.Cells.Item(107, 14).Value = -4400.41174

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(113, 8).Value = 2398
$ws.Cells.Item(113, 10).Value = 2552
$ws.Cells.Item(113, 12).Value = 2552
$ws.Cells.Item(113, 14).Value = -6892

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122, 8).Value = 2607.0527
$ws.Cells.Item(122, 9).Value = 2533.6875
$ws.Cells.Item(122, 11).Value = 7601.0625
$ws.Cells.Item(122, 13).Value = -5151.0625

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 1576.0333
$ws.Cells.Item(132, 9).Value = 1311.1852
$ws.Cells.Item(132, 10).Value = 3959.6667
$ws.Cells.Item(132, 11).Value = 3933.5556
$ws.Cells.Item(132, 12).Value = 11879.0001
$ws.Cells.Item(132, 13).Value = -1403.5556
$ws.Cells.Item(132, 14).Value = -16939.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 2752.6
$ws.Cells.Item(22, 9).Value = 1300
$ws.Cells.Item(22, 11).Value = 1300
$ws.Cells.Item(22, 13).Value = -1005

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(27, 8).Value = 2752.6
$ws.Cells.Item(27, 9).Value = 1300
$ws.Cells.Item(27, 11).Value = 1300
$ws.Cells.Item(27, 13).Value = -1193

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(60, 8).Value = 106000
$ws.Cells.Item(60, 9).Value = 250000
$ws.Cells.Item(60, 11).Value = 250000
$ws.Cells.Item(60, 13).Value = -249491

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(93, 8).Value = 3164.0833
$ws.Cells.Item(93, 10).Value = 2499.75
$ws.Cells.Item(93, 12).Value = 2499.75
$ws.Cells.Item(93, 14).Value = -4995.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 8597.200000000001
$ws.Cells.Item(122, 10).Value = 8749.5
$ws.Cells.Item(122, 12).Value = 26248.5
$ws.Cells.Item(122, 14).Value = -31148.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(126, 8).Value = 3660.7693
$ws.Cells.Item(126, 9).Value = 3590.9167
$ws.Cells.Item(126, 10).Value = 4499
$ws.Cells.Item(126, 11).Value = 10772.7501
$ws.Cells.Item(126, 12).Value = 13497
$ws.Cells.Item(126, 13).Value = -8302.750100000001
$ws.Cells.Item(126, 14).Value = -18437

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 1104.9231
$ws.Cells.Item(136, 9).Value = 736.9
$ws.Cells.Item(136, 11).Value = 2210.7
$ws.Cells.Item(136, 13).Value = 339.3000000000002
